# Menambahkan fitur Rekap Absensi Harian
# Moves the "Hadir" (present) attendance record that was logged against
# 19 Mei 2023 (row 20) back to 12 Mei 2023 (row 13), gives 19 Mei 2023 its
# own fresh in/out record, and clears out the records that had been
# (incorrectly) logged for 21 Mei 2023 (row 22, a Sunday/weekend) and
# 24 Mei 2023 (row 25) - restoring row 22 to a normal "Libur Akhir Pekan"
# weekend row. The daily recap counts at the bottom are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - Jumat, 12 Mei 2023: now has an attendance record.
$ws.Range("B13").Value = "09:43:26"
$ws.Range("D13").Value = "Hadir"
$ws.Range("E13").Value = "14,626 kilometer, TERLAMBAT 2 jam 29 menit"

# Row 20 - Jumat, 19 Mei 2023: replaced with a new in/out record.
$ws.Range("B20").Value = "21:01:29"
$ws.Range("C20").Value = "22:22:48"
$ws.Range("D20").Value = "Hadir"
$ws.Range("E20").Value = "34,163 kilometer, TERLAMBAT 13 jam 47 menit"

# Row 22 - Minggu, 21 Mei 2023: clear the attendance record, back to a
# regular weekend row.
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = ""
$ws.Range("E22").Value = "Libur Akhir Pekan"

# Row 25 - Rabu, 24 Mei 2023: clear the attendance record entirely.
$ws.Range("B25").Value = ""
$ws.Range("D25").Value = ""
$ws.Range("E25").Value = ""

# Recap totals at the bottom (Hadir / Jumlah Keseluruhan) drop from 3 to 2.
$ws.Range("B34").Value = 2
$ws.Range("B37").Value = 2
